$wb = $excel.ActiveWorkbook

# --- Sheet: LP1912 (181 cell changes) ---
$ws = $wb.Worksheets.Item("LP1912")
$ws.Cells.Item(2,1).Value = "Última actualización: 10:39:11"
$ws.Cells.Item(3,1).Value = "Total filas: 172"
$ws.Cells.Item(8,1).Value = "04:44:46"
$ws.Cells.Item(8,3).Value = "15_ABASTO"
$ws.Cells.Item(8,4).Value = 2
$ws.Cells.Item(9,3).Value = "215_EL PELIGRO"
$ws.Cells.Item(10,1).Value = "03:52:04"
$ws.Cells.Item(10,3).Value = "215A_EL PATO"
$ws.Cells.Item(10,4).Value = 54
$ws.Cells.Item(45,1).Value = "06:46:37"
$ws.Cells.Item(45,3).Value = "215A_EL PATO"
$ws.Cells.Item(45,4).Value = 4
$ws.Cells.Item(46,1).Value = "05:16:02"
$ws.Cells.Item(46,3).Value = "17_ROMERO"
$ws.Cells.Item(46,4).Value = 94
$ws.Cells.Item(65,1).Value = "07:12:47"
$ws.Cells.Item(65,3).Value = "27_EL RETIRO"
$ws.Cells.Item(65,4).Value = 24
$ws.Cells.Item(66,1).Value = "06:46:37"
$ws.Cells.Item(66,3).Value = "17X38_ROMERO"
$ws.Cells.Item(66,4).Value = 50
$ws.Cells.Item(75,1).Value = "07:50:33"
$ws.Cells.Item(75,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(75,4).Value = 9
$ws.Cells.Item(76,1).Value = "06:53:56"
$ws.Cells.Item(76,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(76,4).Value = 66
$ws.Cells.Item(93,3).Value = "10_OLMOS"
$ws.Cells.Item(94,3).Value = "215A_EL PATO"
$ws.Cells.Item(99,1).Value = "08:52:13"
$ws.Cells.Item(99,3).Value = "215B_EL PATO"
$ws.Cells.Item(99,4).Value = 8
$ws.Cells.Item(100,1).Value = "08:10:38"
$ws.Cells.Item(100,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(100,4).Value = 50
$ws.Cells.Item(116,1).Value = "09:22:27"
$ws.Cells.Item(116,3).Value = "10_OLMOS"
$ws.Cells.Item(116,4).Value = 4
$ws.Cells.Item(117,1).Value = "08:29:58"
$ws.Cells.Item(117,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(117,4).Value = 57
$ws.Cells.Item(119,1).Value = "08:52:13"
$ws.Cells.Item(119,3).Value = "10_OLMOS"
$ws.Cells.Item(119,4).Value = 37
$ws.Cells.Item(120,1).Value = "08:40:53"
$ws.Cells.Item(120,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(120,4).Value = 49
$ws.Cells.Item(132,1).Value = "08:40:53"
$ws.Cells.Item(132,3).Value = "215C_EL PATO"
$ws.Cells.Item(132,4).Value = 83
$ws.Cells.Item(133,1).Value = "09:22:27"
$ws.Cells.Item(133,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(133,4).Value = 41
$ws.Cells.Item(143,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(144,3).Value = "15_ABASTO"
$ws.Cells.Item(147,1).Value = "10:39:11"
$ws.Cells.Item(147,4).Value = 5
$ws.Cells.Item(148,1).Value = "10:39:11"
$ws.Cells.Item(148,4).Value = 10
$ws.Cells.Item(149,1).Value = "10:39:11"
$ws.Cells.Item(149,4).Value = 12
$ws.Cells.Item(150,1).Value = "10:39:11"
$ws.Cells.Item(150,2).Value = "10:54"
$ws.Cells.Item(150,3).Value = "14_ABASTO"
$ws.Cells.Item(150,4).Value = 15
$ws.Cells.Item(151,1).Value = "09:22:27"
$ws.Cells.Item(151,2).Value = "10:56"
$ws.Cells.Item(151,3).Value = "27_EL RETIRO"
$ws.Cells.Item(151,4).Value = 94
$ws.Cells.Item(152,1).Value = "10:39:11"
$ws.Cells.Item(152,4).Value = 18
$ws.Cells.Item(153,1).Value = "10:39:11"
$ws.Cells.Item(153,2).Value = "10:57"
$ws.Cells.Item(153,3).Value = "14_ABASTO"
$ws.Cells.Item(153,4).Value = 18
$ws.Cells.Item(154,2).Value = "10:57"
$ws.Cells.Item(154,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(154,4).Value = 51
$ws.Cells.Item(155,1).Value = "10:39:11"
$ws.Cells.Item(155,2).Value = "11:01"
$ws.Cells.Item(155,3).Value = "16_SANTA ANA"
$ws.Cells.Item(155,4).Value = 22
$ws.Cells.Item(156,1).Value = "10:06:07"
$ws.Cells.Item(156,2).Value = "11:04"
$ws.Cells.Item(156,3).Value = "17_ROMERO"
$ws.Cells.Item(156,4).Value = 58
$ws.Cells.Item(157,1).Value = "10:39:11"
$ws.Cells.Item(157,2).Value = "11:08"
$ws.Cells.Item(157,3).Value = "225_C ROCA-H SUR"
$ws.Cells.Item(157,4).Value = 29
$ws.Cells.Item(158,2).Value = "11:09"
$ws.Cells.Item(158,3).Value = "17_ROMERO"
$ws.Cells.Item(158,4).Value = 107
$ws.Cells.Item(159,1).Value = "09:22:27"
$ws.Cells.Item(159,2).Value = "11:09"
$ws.Cells.Item(159,3).Value = "14_ABASTO"
$ws.Cells.Item(159,4).Value = 107
$ws.Cells.Item(160,1).Value = "10:39:11"
$ws.Cells.Item(160,2).Value = "11:11"
$ws.Cells.Item(160,3).Value = "17_ROMERO"
$ws.Cells.Item(160,4).Value = 32
$ws.Cells.Item(161,1).Value = "10:39:11"
$ws.Cells.Item(161,2).Value = "11:19"
$ws.Cells.Item(161,3).Value = "215C_EL PATO"
$ws.Cells.Item(161,4).Value = 40
$ws.Cells.Item(162,1).Value = "09:22:27"
$ws.Cells.Item(162,2).Value = "11:20"
$ws.Cells.Item(162,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(162,4).Value = 118
$ws.Cells.Item(163,1).Value = "10:39:11"
$ws.Cells.Item(163,2).Value = "11:21"
$ws.Cells.Item(163,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(163,4).Value = 42
$ws.Cells.Item(164,1).Value = "10:39:11"
$ws.Cells.Item(164,2).Value = "11:21"
$ws.Cells.Item(164,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(164,4).Value = 42
$ws.Cells.Item(165,1).Value = "10:39:11"
$ws.Cells.Item(165,2).Value = "11:30"
$ws.Cells.Item(165,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(165,4).Value = 51
$ws.Cells.Item(166,1).Value = "10:06:07"
$ws.Cells.Item(166,2).Value = "11:30"
$ws.Cells.Item(166,3).Value = "14_ABASTO"
$ws.Cells.Item(166,4).Value = 84
$ws.Cells.Item(166,5).Value = "LP1912"
$ws.Cells.Item(167,1).Value = "10:39:11"
$ws.Cells.Item(167,2).Value = "11:33"
$ws.Cells.Item(167,3).Value = "215A_EL PATO"
$ws.Cells.Item(167,4).Value = 54
$ws.Cells.Item(167,5).Value = "LP1912"
$ws.Cells.Item(168,1).Value = "10:39:11"
$ws.Cells.Item(168,2).Value = "11:41"
$ws.Cells.Item(168,3).Value = "16_SANTA ANA"
$ws.Cells.Item(168,4).Value = 62
$ws.Cells.Item(168,5).Value = "LP1912"
$ws.Cells.Item(169,1).Value = "10:39:11"
$ws.Cells.Item(169,2).Value = "11:45"
$ws.Cells.Item(169,3).Value = "215B_EL PATO"
$ws.Cells.Item(169,4).Value = 66
$ws.Cells.Item(169,5).Value = "LP1912"
$ws.Cells.Item(170,1).Value = "10:39:11"
$ws.Cells.Item(170,2).Value = "11:49"
$ws.Cells.Item(170,3).Value = "15_ABASTO"
$ws.Cells.Item(170,4).Value = 70
$ws.Cells.Item(170,5).Value = "LP1912"
$ws.Cells.Item(171,1).Value = "10:39:11"
$ws.Cells.Item(171,2).Value = "11:51"
$ws.Cells.Item(171,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(171,4).Value = 72
$ws.Cells.Item(171,5).Value = "LP1912"
$ws.Cells.Item(172,1).Value = "10:39:11"
$ws.Cells.Item(172,2).Value = "11:56"
$ws.Cells.Item(172,3).Value = "225_GOMEZ"
$ws.Cells.Item(172,4).Value = 77
$ws.Cells.Item(172,5).Value = "LP1912"
$ws.Cells.Item(173,1).Value = "10:39:11"
$ws.Cells.Item(173,2).Value = "12:04"
$ws.Cells.Item(173,3).Value = "17_ROMERO"
$ws.Cells.Item(173,4).Value = 85
$ws.Cells.Item(173,5).Value = "LP1912"
$ws.Cells.Item(174,1).Value = "10:39:11"
$ws.Cells.Item(174,2).Value = "12:20"
$ws.Cells.Item(174,3).Value = "10_OLMOS"
$ws.Cells.Item(174,4).Value = 101
$ws.Cells.Item(174,5).Value = "LP1912"
$ws.Cells.Item(175,1).Value = "10:39:11"
$ws.Cells.Item(175,2).Value = "12:33"
$ws.Cells.Item(175,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(175,4).Value = 114
$ws.Cells.Item(175,5).Value = "LP1912"
$ws.Cells.Item(176,1).Value = "10:39:11"
$ws.Cells.Item(176,2).Value = "12:34"
$ws.Cells.Item(176,3).Value = "215C_EL PATO"
$ws.Cells.Item(176,4).Value = 115
$ws.Cells.Item(176,5).Value = "LP1912"
$ws.Cells.Item(177,1).Value = "10:39:11"
$ws.Cells.Item(177,2).Value = "12:37"
$ws.Cells.Item(177,3).Value = "27_EL RETIRO"
$ws.Cells.Item(177,4).Value = 118
$ws.Cells.Item(177,5).Value = "LP1912"

# --- Sheet: LP1912-215 (13 cell changes) ---
$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Cells.Item(2,1).Value = "Última actualización: 10:39:11"
$ws.Cells.Item(3,1).Value = "Total filas: 26"
$ws.Cells.Item(28,1).Value = "10:39:11"
$ws.Cells.Item(28,4).Value = 40
$ws.Cells.Item(29,1).Value = "10:39:11"
$ws.Cells.Item(29,4).Value = 54
$ws.Cells.Item(30,1).Value = "10:39:11"
$ws.Cells.Item(30,4).Value = 66
$ws.Cells.Item(31,1).Value = "10:39:11"
$ws.Cells.Item(31,2).Value = "12:34"
$ws.Cells.Item(31,3).Value = "215C_EL PATO"
$ws.Cells.Item(31,4).Value = 115
$ws.Cells.Item(31,5).Value = "LP1912"

# --- Sheet: 6203-6173 (3 cell changes) ---
$ws = $wb.Worksheets.Item("6203-6173")
$ws.Cells.Item(2,1).Value = "Última actualización: 10:39:11"
$ws.Cells.Item(22,1).Value = "10:39:11"
$ws.Cells.Item(22,4).Value = 77
